$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 2.15
$ws.Range("I3").Value = 3.1
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.7
$ws.Range("S3").Value = 1.41
$ws.Range("T3").Value = 2.62
$ws.Range("AG3").Value = 9
$ws.Range("AH3").Value = 15
$ws.Range("AO3").Value = 13
